$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column headers I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold, centered, bordered - same as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I ("I0") and J ("IF") for rows 2-12.
$iValues = @(1, 1, 1, 1, 1, 4, 6, 1, 1, 3, 1)
$jValues = @(6, 4, 6, 5, 7, 7, 8, 5, 6, 4, 2)

for ($r = 0; $r -lt 11; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
